$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the style of the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new "Save" column values for rows 2-9
$saveValues = @(1, 0, 0, 0, 1, 1, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
